# "Generate Report for Handoff"
#
# Two new files finished translation and are ready for handoff:
#   5d37618f-1aaa-4fe0-b22a-b8b37a598912.md
#   881c987d-dcb6-4c24-898f-ed58934bbd6a.md
#
# They are inserted ahead of the existing 9cd57afc-... row (which stays
# "Ready for handoff" too) on all three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

function Add-Hyperlink($ws, $cellRef, $url, $display) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, $null, $null, $display) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2 (7006ebd9...) is untouched.

# Row 3: new file 5d37618f...
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-24-14 04:24:31"
Add-Hyperlink $wsOverview "A3" "https://github.com/OpenLocalizationTest/oltest/blob/197c10903751984469c2934265956ae2cc467201/e2e/5d37618f-1aaa-4fe0-b22a-b8b37a598912.md" "5d37618f-1aaa-4fe0-b22a-b8b37a598912.md"

# Row 4: new file 881c987d...
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-24-14 04:24:31"
Add-Hyperlink $wsOverview "A4" "https://github.com/OpenLocalizationTest/oltest/blob/197c10903751984469c2934265956ae2cc467201/e2e/881c987d-dcb6-4c24-898f-ed58934bbd6a.md" "881c987d-dcb6-4c24-898f-ed58934bbd6a.md"

# Row 5: existing file 9cd57afc..., now pushed down
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-23-14 04:23:28"
Add-Hyperlink $wsOverview "A5" "https://github.com/OpenLocalizationTest/oltest/blob/e20b952efdccbef3a943911969bffe3d3aebde64/e2e/9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.md" "9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": detailed per-locale handoff/handback tracking
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 (7006ebd9...) is untouched.

# Row 3: new file 5d37618f...
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "2016-03-14 04:24:28"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"
Add-Hyperlink $wsZh "A3" "https://github.com/OpenLocalizationTest/oltest/blob/197c10903751984469c2934265956ae2cc467201/e2e/5d37618f-1aaa-4fe0-b22a-b8b37a598912.md" "5d37618f-1aaa-4fe0-b22a-b8b37a598912.md"
Add-Hyperlink $wsZh "B3" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fe1e203e548f7848c3495e576cdbad040eaedd5e/e2e/5d37618f-1aaa-4fe0-b22a-b8b37a598912.md" ".md"
Add-Hyperlink $wsZh "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/80e8162023a723ab6b88fe825cbc49b587e7601e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5d37618f-1aaa-4fe0-b22a-b8b37a598912.07fc38fd80f55ce3046410dd9761d27858a21909.zh-cn.xlf" "5d37618f-1aaa-4fe0-b22a-b8b37a598912.07fc38fd80f55ce3046410dd9761d27858a21909.zh-cn.xlf"

# Row 4: new file 881c987d...
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("E4").Value = "2016-03-14 04:24:28"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "Include"
Add-Hyperlink $wsZh "A4" "https://github.com/OpenLocalizationTest/oltest/blob/197c10903751984469c2934265956ae2cc467201/e2e/881c987d-dcb6-4c24-898f-ed58934bbd6a.md" "881c987d-dcb6-4c24-898f-ed58934bbd6a.md"
Add-Hyperlink $wsZh "B4" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fe1e203e548f7848c3495e576cdbad040eaedd5e/e2e/881c987d-dcb6-4c24-898f-ed58934bbd6a.md" ".md"
Add-Hyperlink $wsZh "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/80e8162023a723ab6b88fe825cbc49b587e7601e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/881c987d-dcb6-4c24-898f-ed58934bbd6a.9fb68ece412a0ba4a56de3c001f34a764c0029af.zh-cn.xlf" "881c987d-dcb6-4c24-898f-ed58934bbd6a.9fb68ece412a0ba4a56de3c001f34a764c0029af.zh-cn.xlf"

# Row 5: existing file 9cd57afc..., now pushed down
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("E5").Value = "2016-03-14 04:23:26"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value = "Include"
Add-Hyperlink $wsZh "A5" "https://github.com/OpenLocalizationTest/oltest/blob/e20b952efdccbef3a943911969bffe3d3aebde64/e2e/9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.md" "9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.md"
Add-Hyperlink $wsZh "B5" "https://github.com/OpenLocalizationTest/oltest/blob/e20b952efdccbef3a943911969bffe3d3aebde64/e2e/9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.md" ".md"
Add-Hyperlink $wsZh "D5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce3f17ce0201569714ade89fb464e6cfc323b41a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.7a717530b67165f4628dfec1daecbc37182f1ad8.zh-cn.xlf" "9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.7a717530b67165f4628dfec1daecbc37182f1ad8.zh-cn.xlf"

# ---------------------------------------------------------------------------
# Sheet "de-de": detailed per-locale handoff/handback tracking
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2 (7006ebd9...) is untouched.

# Row 3: new file 5d37618f...
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "2016-03-14 04:24:31"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"
Add-Hyperlink $wsDe "A3" "https://github.com/OpenLocalizationTest/oltest/blob/197c10903751984469c2934265956ae2cc467201/e2e/5d37618f-1aaa-4fe0-b22a-b8b37a598912.md" "5d37618f-1aaa-4fe0-b22a-b8b37a598912.md"
Add-Hyperlink $wsDe "B3" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fc050dda68d12ac6efb5d351b2ea72cc18f6654b/e2e/5d37618f-1aaa-4fe0-b22a-b8b37a598912.md" ".md"
Add-Hyperlink $wsDe "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a80275b673a20471a6e3477071a0f1da7b17272/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5d37618f-1aaa-4fe0-b22a-b8b37a598912.07fc38fd80f55ce3046410dd9761d27858a21909.de-de.xlf" "5d37618f-1aaa-4fe0-b22a-b8b37a598912.07fc38fd80f55ce3046410dd9761d27858a21909.de-de.xlf"

# Row 4: new file 881c987d...
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("E4").Value = "2016-03-14 04:24:31"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "Include"
Add-Hyperlink $wsDe "A4" "https://github.com/OpenLocalizationTest/oltest/blob/197c10903751984469c2934265956ae2cc467201/e2e/881c987d-dcb6-4c24-898f-ed58934bbd6a.md" "881c987d-dcb6-4c24-898f-ed58934bbd6a.md"
Add-Hyperlink $wsDe "B4" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fc050dda68d12ac6efb5d351b2ea72cc18f6654b/e2e/881c987d-dcb6-4c24-898f-ed58934bbd6a.md" ".md"
Add-Hyperlink $wsDe "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a80275b673a20471a6e3477071a0f1da7b17272/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/881c987d-dcb6-4c24-898f-ed58934bbd6a.9fb68ece412a0ba4a56de3c001f34a764c0029af.de-de.xlf" "881c987d-dcb6-4c24-898f-ed58934bbd6a.9fb68ece412a0ba4a56de3c001f34a764c0029af.de-de.xlf"

# Row 5: existing file 9cd57afc..., now pushed down
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("E5").Value = "2016-03-14 04:23:28"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value = "Include"
Add-Hyperlink $wsDe "A5" "https://github.com/OpenLocalizationTest/oltest/blob/e20b952efdccbef3a943911969bffe3d3aebde64/e2e/9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.md" "9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.md"
Add-Hyperlink $wsDe "B5" "https://github.com/OpenLocalizationTest/oltest/blob/e20b952efdccbef3a943911969bffe3d3aebde64/e2e/9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.md" ".md"
Add-Hyperlink $wsDe "D5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d6ba5c0bf08709287bf9a5bc886785f3e676572c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.7a717530b67165f4628dfec1daecbc37182f1ad8.de-de.xlf" "9cd57afc-7c9d-4d4c-bf8a-1c580689ad6b.7a717530b67165f4628dfec1daecbc37182f1ad8.de-de.xlf"
